$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table lives in A2:D28 (27 daily rows), sorted newest-first.
# A new day's data (2025-12-18) needs to be inserted at the top (row 2),
# pushing every existing row down by one, growing the table to A2:D29.

$lastRow = 28

# Walk bottom-up so we never clobber a source row before reading it.
for ($r = $lastRow; $r -ge 2; $r--) {
    $destRow = $r + 1

    $dateText = $ws.Cells.Item($r, 1).Text
    $bVal = $ws.Cells.Item($r, 2).Value2
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2

    $destA = $ws.Cells.Item($destRow, 1)
    # Force text storage so the yyyy-mm-dd string isn't re-interpreted as a date.
    $destA.NumberFormat = "@"
    $destA.Value = $dateText
    $destA.ClearFormats()

    $ws.Cells.Item($destRow, 2).Value = $bVal
    $ws.Cells.Item($destRow, 3).Value = $cVal
    $ws.Cells.Item($destRow, 4).Value = $dVal
}

# Write the new row for 2025-12-18 into row 2, reusing the same metric values.
$newA = $ws.Cells.Item(2, 1)
$newA.NumberFormat = "@"
$newA.Value = "2025-12-18"
$newA.ClearFormats()

$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
